# Add two new Mac-Addresses (and related device master rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3000166, "Finger Print Scanner 30", "D6-15-AC-80-6B-86", "BS563Q2230814", 165),
    @(3000167, "IRIS Scanner 30",         "6D-58-E2-DF-74-34", "BS563Q2230815", 327),
    @(3000168, "Web Camera 30",           "E2-A8-56-86-15-30", "BS563Q2230816", 736),
    @(3000169, "Document Scanner 30",     "72-E8-B9-FD-63-65", "BS563Q2230817", 801),
    @(3000170, "Printer 30",              "D3-F3-A4-50-AD-12", "BS563Q2230818", 920),
    @(3000171, "Finger Print Scanner 31", "06-16-D0-0B-A6-E4", "BS563Q2230819", 165),
    @(3000172, "IRIS Scanner 31",         "21-78-45-AC-E9-20", "BS563Q2230820", 327),
    @(3000173, "Web Camera 31",           "3C-E8-87-99-DB-FA", "BS563Q2230821", 736),
    @(3000174, "Document Scanner 31",     "BF-55-53-98-40-08", "BS563Q2230822", 801),
    @(3000175, "Printer 31",              "5A-43-36-46-22-EB", "BS563Q2230823", 920)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($r, 9).Value = "superadmin"
    $ws.Cells.Item($r, 10).Value = "now()"
    $ws.Cells.Item($r, 11).Value = "now()"
}

$ws.Range("D145").Select()
